# "Eingabe und erstes einstellen"
#
# The sheet originally listed 4 LED reference numbers (#1..#4), each with its
# own Spannung/Minimalstrom/Maximalstrom row. The table is restructured so
# each reference number has TWO LED current readings (LED1 + LED2) side by
# side, which halves the row count (4 -> 2 data rows) and grows the column
# count (4 -> 6 columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the old #3 / #4 rows (their data gets folded into the #1/#2
#    rows as the new LED2 columns instead).
# ------------------------------------------------------------------
$ws.Rows("4:5").Delete()

# ------------------------------------------------------------------
# 2. Rename the existing current-column headers to be LED1-specific.
# ------------------------------------------------------------------
$ws.Range("C1").Value = "LED1 Mininmalstrom [mA]"
$ws.Range("D1").Value = "LED1 Maximalstrom [mA]"

# ------------------------------------------------------------------
# 3. Old row 3 (#2) used to hold #2's own min/max current (1030/1070).
#    That reading moves to the new LED2 columns (E3/F3); the old #3
#    row's min/max (340/360 -> 320/340 after re-measuring) becomes
#    #2's LED1 reading here instead.
# ------------------------------------------------------------------
$ws.Range("C3").Value = 320
$ws.Range("D3").Value = 340

# ------------------------------------------------------------------
# 4. New LED2 header cells (E1/F1) - copy formatting from a same-role
#    header cell first so they share the existing style/border, then
#    set their text.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("E1").Value = "LED2 Mininmalstrom [mA]"
$ws.Range("F1").Value = "LED2 Maximalstrom [mA]"

# ------------------------------------------------------------------
# 5. New LED2 data cells. F2/F3 share the left/right-border look of
#    columns A-C; E2/E3 get a fresh "separator" look (medium left
#    border only) since it is the first of the two new LED2 columns.
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E2").Value = 1030
$ws.Range("E2").Borders.Item(7).Weight = -4138
$ws.Range("F2").Value = 1070

$ws.Range("E3").Value = 1020
$ws.Range("E3").Borders.Item(7).Weight = -4138
$ws.Range("F3").Value = 1040

# ------------------------------------------------------------------
# 6. Column widths - tuned by hand to match the finished layout.
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 28.5
$ws.Columns("B").ColumnWidth = 12.666666666666666
$ws.Columns("C").ColumnWidth = 25.666666666666668
$ws.Columns("D").ColumnWidth = 24.166666666666668
$ws.Columns("E").ColumnWidth = 25.022135416666668
$ws.Columns("F").ColumnWidth = 23.5

# ------------------------------------------------------------------
# 7. Leave the selection where the author left it.
# ------------------------------------------------------------------
$ws.Range("D6").Select() | Out-Null
